# Add a 2045 ORTP scenario column (column S) to "Proj Attributes and Scenarios"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Proj Attributes and Scenarios")

# Header for the new column
$ws.Range("S9").Value = "ORTP2045"
$ws.Range("S9").Style = $ws.Range("R9").Style

# Fill data rows 10-66 with 1, copying the style used by the rest of the
# scenario columns in each row (style of column R on that row).
for ($r = 10; $r -le 66; $r++) {
    $cell = $ws.Cells.Item($r, 19)
    $cell.Value = 1
    $cell.Style = $ws.Cells.Item($r, 18).Style
}

# Autofit the new column width like the source column.
$ws.Columns.Item(19).AutoFit()

# Extend the worksheet AutoFilter to include the new column.
$ws.Range("A9:S66").AutoFilter()

# Extend the workbook-level _FilterDatabase defined name for this sheet.
$wb.Names.Item("_xlnm._FilterDatabase").RefersToR1C1 = "='Proj Attributes and Scenarios'!R9C1:R66C19"
